$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (SamplesTab): the "Sample ID" query in B3 is rewritten to drop the
# Tumor / Analyte Type columns from the SELECT list (CDS All-studies testcase).
$ws.Cells.Item(3, 2).Value = @"
SELECT
    DISTINCT (smp.sample_id) AS "Sample ID",
    sp.participant_id AS "Participant ID", 
    s.study_name AS "Study Name",
    s.phs_accession AS Accession
FROM 
    df_participant sp
JOIN 
    df_study s ON sp."study.phs_accession" = s.phs_accession
JOIN 
    df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_program p ON p.program_acronym = s."program.program_acronym"
JOIN
    df_file f1 ON f1."sample.sample_id" = smp.sample_id
JOIN
    df_genomic_info gi ON gi."file.file_id" = f1.file_id
WHERE 
  s.phs_accession = 'phs001437' AND d.primary_diagnosis = 'CPNET'
ORDER BY 
    smp.sample_id ASC
LIMIT 100;
"@

# The TsvExcel / WebExcel helper columns (D/E) are no longer populated for
# the SamplesTab and FilesTab rows.
$ws.Cells.Item(3, 4).ClearContents()
$ws.Cells.Item(3, 5).ClearContents()
$ws.Cells.Item(4, 4).ClearContents()
$ws.Cells.Item(4, 5).ClearContents()

# Update the view: scroll up one row and move the active selection to C3.
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 1
[void]$ws.Range("C3").Select()
